$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above row 4, pushing the existing rows 4-6
# (Pol0_45_90_135 / Img1_C2.tif / Img2_C2.tif) down to rows 6-8.
$ws.Rows("4:5").Insert(-4121)

# The newly inserted rows inherit the formatting of the row above; fully
# clear them (content + formatting) so they serialize as blank rows.
$ws.Range("A4:D5").Clear()

# Update the selection to match the target state (active cell A6, selection A6:A8)
$ws.Range("A6:A8").Select()
